$d = $word.ActiveDocument

function Find-ParagraphByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# Change 2: move <w:lastRenderedPageBreak/> from the "9.53am 16/05/2017"
# paragraph to the start of the "9.01pm 13/05/2017" paragraph.
# ------------------------------------------------------------------
$pPageBreakSource = Find-ParagraphByText("9.53am 16/05/2017")
$xmlNoBreak = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>9.53am 16/05/2017 – Updated Gantt Chart and made it more comprehensive.</w:t></w:r></w:p>'
$pPageBreakSource.Range.InsertXML($xmlNoBreak)

$pPageBreakTarget = Find-ParagraphByText("9.01pm 13/05/2017")
$xmlWithBreak = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>9.01pm 13/05/2017 – Updated portfolio and improved some sections</w:t></w:r></w:p>'
$pPageBreakTarget.Range.InsertXML($xmlWithBreak)

# ------------------------------------------------------------------
# Change 3: split the "8.05pm 22/05/2017" paragraph into two runs with
# the "_GoBack" bookmark moved in between them.
# ------------------------------------------------------------------
$pDataDict = Find-ParagraphByText("8.05pm 22/05/2017")
$xmlSplit = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">8.05pm 22/05/2017 – Updated Data Dictionary. Added testing information and feedback. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Added MIT License.</w:t></w:r></w:p>'
$pDataDict.Range.InsertXML($xmlSplit)

# ------------------------------------------------------------------
# Change 4: remove the old "_GoBack" bookmark from the final paragraph
# (it previously trailed the "10.21pm 11/06/2017" paragraph text).
# ------------------------------------------------------------------
$pLast = Find-ParagraphByText("10.21pm 11/06/2017")
$xmlLast = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>10.21pm 11/06/2017 – Collected another survey response which had found 12 bugs. These bugs were submitted on GitHub. Fixed 2 major bugs including one of the algorithms not taking into account the landing height, and Microsoft Excel not closing and staying open according to Task Manager in the background.</w:t></w:r></w:p>'
$pLast.Range.InsertXML($xmlLast)

# ------------------------------------------------------------------
# Change 1: insert two new logbook entries before the
# "12.23pm 27/04/2017" paragraph.
# ------------------------------------------------------------------
$pFeedback = Find-ParagraphByText("12.23pm 27/04/2017")
$newEntry1 = "9.26pm 20/04/2017 – Fixed a bug where the program did not accept any other acceleration value apart from 9.8ms^-2. Also fixed a bug where the angle was given in radians, not in degrees.`r"
$newEntry2 = "10.59pm 23/04/2017 – Fixed a bug where “300” was hard coded into the program instead of the variable “range”.`r"
$pFeedback.Range.InsertBefore($newEntry1 + $newEntry2)

Write-Output "done"
